$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "JUASu713"
$ws.Range("B2").Value = 23111038
$ws.Range("C2").Value = "uigpugd85"
$ws.Range("D2").Value = "U8&Nv%t4"
$ws.Range("F2").Value = "KbRZTWlq"
$ws.Range("G2").Value = "Iqlp"

# Row 3 updates
$ws.Range("A3").Value = "WMUmm578"
$ws.Range("B3").Value = 23111037
$ws.Range("C3").Value = "uidlbnc89"
$ws.Range("D3").Value = "D#N8p4d!"
$ws.Range("F3").Value = "DdpGcIWk"
$ws.Range("G3").Value = "GEff"
